# "ajustes en el users para el add sucursales_id"
#
# Moves the hidden "_GoBack" bookmark from the start of the
# "...privilegios de super admin." bullet to a point inside the
# "...asignados a su sucursal..." bullet, right after the word
# "sucursal" (splitting that run in two) and highlights the first
# half ("Podrá registrar usuarios asignados a su sucursal") in cyan.

$d = $word.ActiveDocument

# Locate the split point: end of "Podrá registrar usuarios asignados
# a su sucursal" inside the target bullet paragraph.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Podrá registrar usuarios asignados a su sucursal", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitStart = $anchor.Start
$splitEnd = $anchor.End

# Move the "_GoBack" bookmark from its old location to the split point.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$bmRange = $d.Range($splitEnd, $splitEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Highlight just the first half of the (now split) run in cyan.
$firstHalf = $d.Range($splitStart, $splitEnd)
$firstHalf.Font.HighlightColorIndex = 3
